$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1000
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H43").Value = 4312.8667
$ws.Range("I43").Value = 925
$ws.Range("J43").Value = 5544.8184
$ws.Range("K43").Value = 925
$ws.Range("L43").Value = 5544.8184
$ws.Range("M43").Value = -856
$ws.Range("N43").Value = -5682.8184
$ws.Range("H69").Value = 1535.7407
$ws.Range("J69").Value = 1506.3462
$ws.Range("L69").Value = 4519.0386
$ws.Range("N69").Value = -6267.0386
$ws.Range("H70").Value = 972.9655
$ws.Range("I70").Value = 991.73914
$ws.Range("J70").Value = 901
$ws.Range("K70").Value = 2975.21742
$ws.Range("L70").Value = 2703
$ws.Range("M70").Value = -2705.21742
$ws.Range("N70").Value = -3243
$ws.Range("H72").Value = 1535.7407
$ws.Range("J72").Value = 1506.3462
$ws.Range("L72").Value = 13557.1158
$ws.Range("N72").Value = -22293.1158
$ws.Range("H73").Value = 972.9655
$ws.Range("I73").Value = 991.73914
$ws.Range("J73").Value = 901
$ws.Range("K73").Value = 2975.21742
$ws.Range("L73").Value = 2703
$ws.Range("M73").Value = -2039.21742
$ws.Range("N73").Value = -4575
$ws.Range("H76").Value = 3706478.2
$ws.Range("I76").Value = 2929.5
$ws.Range("J76").Value = 11113575
$ws.Range("K76").Value = 2929.5
$ws.Range("L76").Value = 11113575
$ws.Range("M76").Value = -2614.5
$ws.Range("N76").Value = -11114205
$ws.Range("H79").Value = 3706478.2
$ws.Range("I79").Value = 2929.5
$ws.Range("J79").Value = 11113575
$ws.Range("K79").Value = 2929.5
$ws.Range("L79").Value = 11113575
$ws.Range("M79").Value = -1837.5
$ws.Range("N79").Value = -11115759
$ws.Range("H98").Value = 867.63635
$ws.Range("I98").Value = 905.5
$ws.Range("J98").Value = 766.6667
$ws.Range("K98").Value = 905.5
$ws.Range("L98").Value = 766.6667
$ws.Range("M98").Value = 592.5
$ws.Range("N98").Value = -3762.6667
$ws.Range("H112").Value = 10733391
$ws.Range("J112").Value = 5556714.5
$ws.Range("L112").Value = 16670143.5
$ws.Range("N112").Value = -16672359.5
$ws.Range("H122").Value = 867.63635
$ws.Range("I122").Value = 905.5
$ws.Range("J122").Value = 766.6667
$ws.Range("K122").Value = 2716.5
$ws.Range("L122").Value = 2300.0001
$ws.Range("M122").Value = -266.5
$ws.Range("N122").Value = -7200.0001
$ws.Range("H123").Value = 29497.5
$ws.Range("J123").Value = 29497.5
$ws.Range("L123").Value = 29497.5
$ws.Range("N123").Value = -39297.5
$ws.Range("H129").Value = 205015.75
$ws.Range("J129").Value = 223210.48
$ws.Range("L129").Value = 669631.4400000001
$ws.Range("N129").Value = -679631.4400000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1925.8
$ws.Range("I2").Value = 1823.9166
$ws.Range("K2").Value = 1823.9166
$ws.Range("M2").Value = -1710.9166
$ws.Range("H3").Value = 362
$ws.Range("I3").Value = 836.6667
$ws.Range("J3").Value = 6
$ws.Range("K3").Value = 836.6667
$ws.Range("L3").Value = 6
$ws.Range("M3").Value = -721.6667
$ws.Range("N3").Value = -236
$ws.Range("H32").Value = 14901.627
$ws.Range("I32").Value = 10122.796
$ws.Range("J32").Value = 34752.152
$ws.Range("K32").Value = 10122.796
$ws.Range("L32").Value = 34752.152
$ws.Range("M32").Value = -9835.796
$ws.Range("N32").Value = -35326.152
$ws.Range("H45").Value = 4092.8076
$ws.Range("I45").Value = 4038.0557
$ws.Range("J45").Value = 4216
$ws.Range("K45").Value = 4038.0557
$ws.Range("L45").Value = 4216
$ws.Range("M45").Value = -3661.0557
$ws.Range("N45").Value = -4970
$ws.Range("H116").Value = 1925.8
$ws.Range("I116").Value = 1823.9166
$ws.Range("K116").Value = 1823.9166
$ws.Range("M116").Value = 470.0834
$ws.Range("H132").Value = 8855.25
$ws.Range("I132").Value = 1572.55
$ws.Range("J132").Value = 45268.75
$ws.Range("K132").Value = 4717.65
$ws.Range("L132").Value = 135806.25
$ws.Range("M132").Value = -2187.65
$ws.Range("N132").Value = -140866.25
$ws.Range("H139").Value = 50128.8
$ws.Range("J139").Value = 50128.8
$ws.Range("L139").Value = 50128.8
$ws.Range("N139").Value = -60408.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1925.8
$ws.Range("I3").Value = 1823.9166
$ws.Range("K3").Value = 1823.9166
$ws.Range("M3").Value = -1709.9166
$ws.Range("H99").Value = 1605.25
$ws.Range("I99").Value = 1805
$ws.Range("J99").Value = 1405.5
$ws.Range("K99").Value = 1805
$ws.Range("L99").Value = 1405.5
$ws.Range("M99").Value = -307
$ws.Range("N99").Value = -4401.5
$ws.Range("H107").Value = 3127.5833
$ws.Range("I107").Value = 2401.1177
$ws.Range("J107").Value = 4891.857
$ws.Range("K107").Value = 2401.1177
$ws.Range("L107").Value = 4891.857
$ws.Range("M107").Value = -481.1176999999998
$ws.Range("N107").Value = -8731.857
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 174.54546
$ws.Range("I22").Value = 174.54546
$ws.Range("K22").Value = 174.54546
$ws.Range("M22").Value = 175.45454
$ws.Range("H31").Value = 6552.159
$ws.Range("J31").Value = 8043.1665
$ws.Range("L31").Value = 8043.1665
$ws.Range("N31").Value = -8633.166499999999
$ws.Range("H34").Value = 6552.159
$ws.Range("J34").Value = 8043.1665
$ws.Range("L34").Value = 8043.1665
$ws.Range("N34").Value = -8447.166499999999
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H70").Value = 12000
$ws.Range("J70").Value = 12000
$ws.Range("L70").Value = 12000
$ws.Range("N70").Value = -12630
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H73").Value = 12000
$ws.Range("J73").Value = 12000
$ws.Range("L73").Value = 12000
$ws.Range("N73").Value = -14184
$ws.Range("H132").Value = 2326.5
$ws.Range("I132").Value = 1708.3611
$ws.Range("K132").Value = 5125.0833
$ws.Range("M132").Value = -2595.0833
$ws.Range("H141").Value = 25240
$ws.Range("J141").Value = 26764
$ws.Range("L141").Value = 26764
$ws.Range("N141").Value = -37124
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3333407.8
$ws.Range("I4").Value = 95.57143000000001
$ws.Range("J4").Value = 15000000
$ws.Range("K4").Value = 286.71429
$ws.Range("L4").Value = 45000000
$ws.Range("M4").Value = -174.71429
$ws.Range("N4").Value = -45000224
$ws.Range("H59").Value = 3000
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H131").Value = 731.09
$ws.Range("J131").Value = 731.09
$ws.Range("L131").Value = 2193.27
$ws.Range("N131").Value = -12273.27
$ws.Range("H134").Value = 4708.3335
$ws.Range("I134").Value = 2902.5
$ws.Range("J134").Value = 7116.1113
$ws.Range("K134").Value = 8707.5
$ws.Range("L134").Value = 21348.3339
$ws.Range("M134").Value = -3637.5
$ws.Range("N134").Value = -31488.3339
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 11302.8
$ws.Range("I122").Value = 9628.5
$ws.Range("K122").Value = 28885.5
$ws.Range("M122").Value = -26435.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2334.16
$ws.Range("I100").Value = 1805.4
$ws.Range("J100").Value = 2686.6667
$ws.Range("K100").Value = 1805.4
$ws.Range("L100").Value = 2686.6667
$ws.Range("M100").Value = -1264.4
$ws.Range("N100").Value = -3768.6667
$ws.Range("H122").Value = 1249352.8
$ws.Range("I122").Value = 1332176.4
$ws.Range("K122").Value = 3996529.2
$ws.Range("M122").Value = -3994079.2
$ws.Range("H132").Value = 4750
$ws.Range("I132").Value = 1333.3334
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 4000.0002
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -1470.0002
$ws.Range("N132").Value = -50060
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 5000
$ws.Range("J86").Value = 5000
$ws.Range("L86").Value = 5000
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 5000
$ws.Range("J89").Value = 5000
$ws.Range("L89").Value = 25000
$ws.Range("N89").Value = -36232
$ws.Range("H132").Value = 3666.3333
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 4749.5
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 14248.5
$ws.Range("M132").Value = -1970
$ws.Range("N132").Value = -19308.5
